$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 10762.9
$ws.Range("I9").Value = 12766.25
$ws.Range("K9").Value = 12766.25
$ws.Range("M9").Value = -12597.25

$ws.Range("H17").Value = 2022208.4
$ws.Range("J17").Value = 2022208.4
$ws.Range("L17").Value = 6066625.199999999
$ws.Range("N17").Value = -6066961.199999999

$ws.Range("H92").Value = 1581.8334
$ws.Range("I92").Value = 1454
$ws.Range("K92").Value = 1454
$ws.Range("M92").Value = -206

$ws.Range("H116").Value = 6311.846
$ws.Range("I116").Value = 6381.875
$ws.Range("J116").Value = 6199.8
$ws.Range("K116").Value = 6381.875
$ws.Range("L116").Value = 6199.8
$ws.Range("M116").Value = -2939.875
$ws.Range("N116").Value = -13083.8

$ws.Range("H138").Value = 8552151
$ws.Range("I138").Value = 1741
$ws.Range("J138").Value = 13340380
$ws.Range("K138").Value = 5223
$ws.Range("L138").Value = 40021140
$ws.Range("M138").Value = -83
$ws.Range("N138").Value = -40031420

$ws.Range("H141").Value = 3202.375
$ws.Range("I141").Value = 3202.375
$ws.Range("K141").Value = 9607.125
$ws.Range("M141").Value = -4427.125

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2363.0833
$ws.Range("I2").Value = 2557
$ws.Range("K2").Value = 2557
$ws.Range("M2").Value = -2444

$ws.Range("H32").Value = 31261324
$ws.Range("I32").Value = 45461180
$ws.Range("J32").Value = 21636.4
$ws.Range("K32").Value = 45461180
$ws.Range("L32").Value = 21636.4
$ws.Range("M32").Value = -45460893
$ws.Range("N32").Value = -22210.4

$ws.Range("H45").Value = 2436.7
$ws.Range("I45").Value = 1733
$ws.Range("J45").Value = 2738.2856
$ws.Range("K45").Value = 1733
$ws.Range("L45").Value = 2738.2856
$ws.Range("M45").Value = -1356
$ws.Range("N45").Value = -3492.2856

$ws.Range("H102").Value = 2249.6365
$ws.Range("I102").Value = 1415
$ws.Range("K102").Value = 1415
$ws.Range("M102").Value = 207

$ws.Range("H116").Value = 2363.0833
$ws.Range("I116").Value = 2557
$ws.Range("K116").Value = 2557
$ws.Range("M116").Value = -263

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2363.0833
$ws.Range("I3").Value = 2557
$ws.Range("K3").Value = 2557
$ws.Range("M3").Value = -2443

$ws.Range("H20").Value = 3463.2727
$ws.Range("I20").Value = 3637.375
$ws.Range("K20").Value = 3637.375
$ws.Range("M20").Value = -3390.375

$ws.Range("H86").Value = 12634.046
$ws.Range("I86").Value = 3079.7856
$ws.Range("K86").Value = 3079.7856
$ws.Range("M86").Value = -1956.7856

$ws.Range("H89").Value = 12634.046
$ws.Range("I89").Value = 3079.7856
$ws.Range("K89").Value = 15398.928
$ws.Range("M89").Value = -9782.928

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1219.875
$ws.Range("I16").Value = 1115.7142
$ws.Range("K16").Value = 1115.7142
$ws.Range("M16").Value = -828.7141999999999

$ws.Range("H86").Value = 3338.5
$ws.Range("I86").Value = 3194.9167
$ws.Range("J86").Value = 4200
$ws.Range("K86").Value = 3194.9167
$ws.Range("L86").Value = 4200
$ws.Range("M86").Value = -2071.9167
$ws.Range("N86").Value = -6446

$ws.Range("H89").Value = 3338.5
$ws.Range("I89").Value = 3194.9167
$ws.Range("J89").Value = 4200
$ws.Range("K89").Value = 15974.5835
$ws.Range("L89").Value = 21000
$ws.Range("M89").Value = -10358.5835
$ws.Range("N89").Value = -32232

$ws.Range("H113").Value = 1219.875
$ws.Range("I113").Value = 1115.7142
$ws.Range("K113").Value = 1115.7142
$ws.Range("M113").Value = 1054.2858

$ws.Range("H134").Value = 1679.2941
$ws.Range("I134").Value = 1632.0714
$ws.Range("J134").Value = 1899.6666
$ws.Range("K134").Value = 4896.2142
$ws.Range("L134").Value = 5698.9998
$ws.Range("M134").Value = -2361.2142
$ws.Range("N134").Value = -10768.9998

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("M14").ClearContents()
$ws.Range("H14").Value = 0
$ws.Range("I14").Value = 0
$ws.Range("K14").Value = 0

$ws.Range("M88").ClearContents()
$ws.Range("H88").Value = 20000
$ws.Range("I88").Value = 0
$ws.Range("J88").Value = 20000
$ws.Range("K88").Value = 0
$ws.Range("L88").Value = 60000
$ws.Range("N88").Value = -60856

$ws.Range("M91").ClearContents()
$ws.Range("H91").Value = 20000
$ws.Range("I91").Value = 0
$ws.Range("J91").Value = 20000
$ws.Range("K91").Value = 0
$ws.Range("L91").Value = 60000
$ws.Range("N91").Value = -62964

$ws.Range("H107").Value = 1012
$ws.Range("J107").Value = 1800
$ws.Range("L107").Value = 5400
$ws.Range("N107").Value = -9240

$ws.Range("H131").Value = 32757.889
$ws.Range("J131").Value = 5855.364
$ws.Range("L131").Value = 17566.092
$ws.Range("N131").Value = -27646.092

$ws.Range("H132").Value = 2383932
$ws.Range("J132").Value = 3512069.5
$ws.Range("L132").Value = 31608625.5
$ws.Range("N132").Value = -31613685.5

$ws.Range("H133").Value = 13602.2
$ws.Range("J133").Value = 19603.8
$ws.Range("L133").Value = 58811.39999999999
$ws.Range("N133").Value = -68931.39999999999

$ws.Range("H134").Value = 8305.6
$ws.Range("I134").Value = 2277.923
$ws.Range("K134").Value = 6833.768999999999
$ws.Range("M134").Value = -1763.768999999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 7836
$ws.Range("J70").Value = 10000
$ws.Range("L70").Value = 10000
$ws.Range("N70").Value = -10540

$ws.Range("H73").Value = 7836
$ws.Range("J73").Value = 10000
$ws.Range("L73").Value = 10000
$ws.Range("N73").Value = -11872

$ws.Range("H96").Value = 17604.5
$ws.Range("J96").Value = 16806
$ws.Range("L96").Value = 16806
$ws.Range("N96").Value = -22298

$ws.Range("H97").Value = 2634.1428
$ws.Range("I97").Value = 1151.3334
$ws.Range("J97").Value = 3746.25
$ws.Range("K97").Value = 1151.3334
$ws.Range("L97").Value = 3746.25
$ws.Range("M97").Value = -655.3334
$ws.Range("N97").Value = -4738.25

$ws.Range("H122").Value = 1994.2
$ws.Range("I122").Value = 1772.9231
$ws.Range("K122").Value = 5318.7693
$ws.Range("M122").Value = -2868.7693

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 527.44446
$ws.Range("I16").Value = 527.44446
$ws.Range("K16").Value = 527.44446
$ws.Range("M16").Value = -357.44446

$ws.Range("H22").Value = 2698.3225
$ws.Range("I22").Value = 1950.0834
$ws.Range("J22").Value = 3170.8948
$ws.Range("K22").Value = 1950.0834
$ws.Range("L22").Value = 3170.8948
$ws.Range("M22").Value = -1655.0834
$ws.Range("N22").Value = -3760.8948

$ws.Range("H27").Value = 2698.3225
$ws.Range("I27").Value = 1950.0834
$ws.Range("J27").Value = 3170.8948
$ws.Range("K27").Value = 1950.0834
$ws.Range("L27").Value = 3170.8948
$ws.Range("M27").Value = -1843.0834
$ws.Range("N27").Value = -3384.8948

$ws.Range("H46").Value = 1295.7142
$ws.Range("I46").Value = 648
$ws.Range("K46").Value = 648
$ws.Range("M46").Value = -460

$ws.Range("H122").Value = 3576.9119
$ws.Range("I122").Value = 3088.64
$ws.Range("J122").Value = 4933.222
$ws.Range("K122").Value = 9265.92
$ws.Range("L122").Value = 14799.666
$ws.Range("M122").Value = -6815.92
$ws.Range("N122").Value = -19699.666

$ws.Range("H131").Value = 88578
$ws.Range("J131").Value = 88578
$ws.Range("L131").Value = 88578
$ws.Range("N131").Value = -98658

$ws.Range("H136").Value = 1911.2222
$ws.Range("I136").Value = 1935.14
$ws.Range("J136").Value = 1612.25
$ws.Range("K136").Value = 5805.42
$ws.Range("L136").Value = 4836.75
$ws.Range("M136").Value = -3255.42
$ws.Range("N136").Value = -9936.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 818.25
$ws.Range("I81").Value = 595
$ws.Range("K81").Value = 1190
$ws.Range("M81").Value = -129

$ws.Range("H84").Value = 818.25
$ws.Range("I84").Value = 595
$ws.Range("K84").Value = 5950
$ws.Range("M84").Value = -646

$ws.Range("H95").Value = 29997.8
$ws.Range("J95").Value = 29997.8
$ws.Range("L95").Value = 29997.8
$ws.Range("N95").Value = -35489.8

$ws.Range("H132").Value = 4527.298
$ws.Range("I132").Value = 4848.237
$ws.Range("J132").Value = 3172.2222
$ws.Range("K132").Value = 14544.711
$ws.Range("L132").Value = 9516.6666
$ws.Range("M132").Value = -12014.711
$ws.Range("N132").Value = -14576.6666
